# Contact Manager.pptx — "Add more info to the slides"
#
# Summary of the edit (reconstructed from the target OOXML diff):
#  - Slide 2 ("The DB+J's"): add role annotations to the team member bullets.
#  - Slide 3 ("Stacks on Stacks"): tweak the last two stack bullets.
#  - Insert a new slide "Tools" after "Use Case Diagram".
#  - Slide "How did we do it?": add bullet content describing the build steps.
#  - Slide "What went well?": trim the title and add bullet content.
#  - Insert a new slide "WhaT Went Wrong" after "What went well?".
#  - "DEMo" and "Any Questions?" slides shift later in the deck but are
#    otherwise untouched.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------------
# Helper: replace the first occurrence of $old with $new inside a shape's
# text, using a Characters() sub-range so the surrounding run formatting
# (rPr) of untouched runs is preserved.
# ---------------------------------------------------------------------------
function Replace-InShapeText {
    param($shape, [string]$old, [string]$new)
    $tr = $shape.TextFrame.TextRange
    $full = $tr.Text
    $idx = $full.IndexOf($old)
    if ($idx -ge 0) {
        $sub = $tr.Characters($idx + 1, $old.Length)
        $sub.Text = $new
    }
}

# ---------------------------------------------------------------------------
# Slide 2 — "The DB+J's": add role labels to the back-end team bullets.
# ---------------------------------------------------------------------------
$s2 = $p.Slides.Item(2)
$body2 = $s2.Shapes.Item(2)

Replace-InShapeText $body2 "Fochezato" "Fochezato – Database"
Replace-InShapeText $body2 "Blake Robertson" "Blake Robertson - API"
Replace-InShapeText $body2 "Jacob Thomas" "Jacob Thomas - API"

# ---------------------------------------------------------------------------
# Slide 3 — "Stacks on Stacks": update the last two bullets.
# ---------------------------------------------------------------------------
$s3 = $p.Slides.Item(3)
$body3 = $s3.Shapes.Item(2)

Replace-InShapeText $body3 "React.js (API)" "HTML & CSS (Front-End)"
Replace-InShapeText $body3 "Node.js (server)" "Node.js ("

# ---------------------------------------------------------------------------
# New slide — "Tools" (inserted after "Use Case Diagram", i.e. index 6).
# ---------------------------------------------------------------------------
$toolsSlide = $p.Slides.Add(6, 16)
$toolsSlide.Shapes.Title.TextFrame.TextRange.Text = "Tools"
$toolsBody = $toolsSlide.Shapes.Item(2)
$toolsBody.TextFrame.TextRange.Text = "POSTMAN – testing json requests`rJquery`rAWS - `r`r`r"

# ---------------------------------------------------------------------------
# "How did we do it?" (now at index 7) — add the build-steps bullets.
# ---------------------------------------------------------------------------
$howSlide = $p.Slides.Item(7)
$howBody = $howSlide.Shapes.Item(2)
$howBody.TextFrame.TextRange.Text = "Begin with the models for the user and contact in the database and the HTML skeleton`rCreated basic routes for the website`rMade basic json get and post request for the API`rTesting the API with the database`rTesting the API with the front-end`rMerge it all together"

# ---------------------------------------------------------------------------
# "What went well? ..." (now at index 8) — trim title, add bullets.
# ---------------------------------------------------------------------------
$wellSlide = $p.Slides.Item(8)
Replace-InShapeText $wellSlide.Shapes.Title "What went well? What didn't? Improvements for the future." "What went well? Improvements for the future."
$wellBody = $wellSlide.Shapes.Item(2)
$wellBody.TextFrame.TextRange.Text = "Setup communication very early `rHad assigned roles`rSemi-Daily Meetups`rWorking demo before presentation`rHave md5 hash for passwords`r"

# ---------------------------------------------------------------------------
# New slide — "WhaT Went Wrong" (inserted after "What went well?", index 9).
# ---------------------------------------------------------------------------
$wrongSlide = $p.Slides.Add(9, 16)
$wrongSlide.Shapes.Title.TextFrame.TextRange.Text = "WhaT Went Wrong"
$wrongBody = $wrongSlide.Shapes.Item(2)
$wrongBody.TextFrame.TextRange.Text = "No members had any experience in web development`rMERN stack difficult stack for beginners in web development`rShould of spent less time on researching`rNo organized tasks`rDifficulties with connecting API with database and front-end`rLittle to no experience with using Github for most members`r"
